# BOT; UPDATE DATA
# Appends one more day of data to the "相談件数" (consultation counts) sheet.
#
# Layout on sheet 1 (相談件数):
#   - Row 90 was a blank placeholder row (styles already in place: A=3,B/C=7,D/E=8)
#   - Row 91 held the footnote text ("*4/8...") in column B
#
# This update:
#   1. Fills the previously-blank row 90 with the new day's figures (2020-04-24 / serial 43945)
#   2. Inserts a fresh blank placeholder row at 91 (shifting the footnote row down to 92),
#      and stamps it with the next day's date (2020-04-25 / serial 43946) so it is ready
#      for tomorrow's update, matching the existing pattern used throughout the sheet.
#   3. Updates the print area to keep the footnote row included.
#   4. Moves the selection to reflect the new bottom of the data entry area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")
$ws.Activate()

# --- 1) Insert a new placeholder row above the footnote row ------------------
# This shifts the footnote ("B91" containing the shared string) down to row 92,
# and the newly inserted row 91 inherits row 90's formatting (same as Excel's
# native Insert behaviour when inserting a row that copies formats from above).
$ws.Rows.Item(91).Insert()

# --- 2) Fill in today's data into row 90 --------------------------------------
$ws.Range("A90").Value = 43945
$ws.Range("B90").Value = 609
$ws.Range("C90").Value = 28990
$ws.Range("D90").Value = 104
$ws.Range("E90").Value = 6458

# --- 3) Stamp the new blank placeholder row (91) with tomorrow's date --------
$ws.Range("A91").Value = 43946

# --- 4) Extend the print area to still include the (now shifted) footnote ---
$wb.Names.Item(1).RefersTo = "=相談件数!`$A`$1:`$E`$94"

# --- 5) Update the recorded selection to the new last data cell -------------
$null = $ws.Range("E93").Select()
